$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -1
$ws.Range("B1").Value = 2.880195140838623
$ws.Range("C1").Value = 2.593844652175903
$ws.Range("D1").Value = 2.899602890014648
$ws.Range("E1").Value = -1
